$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the STUDENT_ID column header to REGISTRATION NO
$ws.Range("A2").Value = "REGISTRATION NO"

# Update the selected/active cell shown when the sheet is opened
$ws.Range("C4").Select()
